$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 575.5714
$ws.Range("I39").Value = 504.83334
$ws.Range("K39").Value = 1514.50002
$ws.Range("M39").Value = -1218.50002
$ws.Range("I40").Value = 10000
$ws.Range("K40").Value = 10000
$ws.Range("M40").Value = -9825
$ws.Range("H70").Value = 1906.1333
$ws.Range("I70").Value = 1385.2858
$ws.Range("K70").Value = 4155.857400000001
$ws.Range("M70").Value = -3885.857400000001
$ws.Range("H73").Value = 1906.1333
$ws.Range("I73").Value = 1385.2858
$ws.Range("K73").Value = 4155.857400000001
$ws.Range("M73").Value = -3219.857400000001
$ws.Range("H92").Value = 1592.7273
$ws.Range("I92").Value = 1592.7273
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1592.7273
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -344.7273
$ws.Range("N92").ClearContents()
$ws.Range("H137").Value = 605
$ws.Range("I137").Value = 605
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 1815
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 735
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 2577.4285
$ws.Range("I138").Value = 2097.7144
$ws.Range("J138").Value = 3057.1428
$ws.Range("K138").Value = 6293.1432
$ws.Range("L138").Value = 9171.428400000001
$ws.Range("M138").Value = -1153.1432
$ws.Range("N138").Value = -19451.4284
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1306.125
$ws.Range("I5").Value = 84.75
$ws.Range("J5").Value = 2527.5
$ws.Range("K5").Value = 84.75
$ws.Range("L5").Value = 2527.5
$ws.Range("M5").Value = 27.25
$ws.Range("N5").Value = -2751.5
$ws.Range("H32").Value = 2693.5908
$ws.Range("I32").Value = 2356.4285
$ws.Range("K32").Value = 2356.4285
$ws.Range("M32").Value = -2069.4285
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1306.125
$ws.Range("I4").Value = 84.75
$ws.Range("J4").Value = 2527.5
$ws.Range("K4").Value = 84.75
$ws.Range("L4").Value = 2527.5
$ws.Range("M4").Value = 30.25
$ws.Range("N4").Value = -2757.5
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H99").Value = 3048.6667
$ws.Range("I99").Value = 3048.6667
$ws.Range("K99").Value = 3048.6667
$ws.Range("M99").Value = -1550.6667
$ws.Range("H109").Value = 33333
$ws.Range("J109").Value = 33333
$ws.Range("L109").Value = 33333
$ws.Range("N109").Value = -36107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H115").Value = 44245
$ws.Range("I115").Value = 44245
$ws.Range("K115").Value = 44245
$ws.Range("M115").Value = -43070
$ws.Range("H134").Value = 2160.7778
$ws.Range("I134").Value = 2054.125
$ws.Range("K134").Value = 6162.375
$ws.Range("M134").Value = -3627.375
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 153.53847
$ws.Range("I38").Value = 52.625
$ws.Range("J38").Value = 315
$ws.Range("K38").Value = 157.875
$ws.Range("L38").Value = 945
$ws.Range("M38").Value = 189.125
$ws.Range("N38").Value = -1639
$ws.Range("H75").Value = 2749.5
$ws.Range("I75").Value = 499
$ws.Range("J75").Value = 5000
$ws.Range("K75").Value = 1497
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = -499
$ws.Range("N75").Value = -16996
$ws.Range("H78").Value = 2749.5
$ws.Range("I78").Value = 499
$ws.Range("J78").Value = 5000
$ws.Range("K78").Value = 4491
$ws.Range("L78").Value = 45000
$ws.Range("M78").Value = 501
$ws.Range("N78").Value = -54984
$ws.Range("H97").Value = 785.1429000000001
$ws.Range("I97").Value = 947.5
$ws.Range("J97").Value = 720.2
$ws.Range("K97").Value = 2842.5
$ws.Range("L97").Value = 2160.6
$ws.Range("M97").Value = -2346.5
$ws.Range("N97").Value = -3152.6
$ws.Range("H131").Value = 5000
$ws.Range("J131").Value = 5000
$ws.Range("L131").Value = 15000
$ws.Range("N131").Value = -25080
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H70").Value = 7000
$ws.Range("J70").Value = 7000
$ws.Range("L70").Value = 7000
$ws.Range("N70").Value = -7540
$ws.Range("H73").Value = 7000
$ws.Range("J73").Value = 7000
$ws.Range("L73").Value = 7000
$ws.Range("N73").Value = -8872
$ws.Range("H97").Value = 1730
$ws.Range("I97").Value = 1593.3334
$ws.Range("J97").Value = 2003.3334
$ws.Range("K97").Value = 1593.3334
$ws.Range("L97").Value = 2003.3334
$ws.Range("M97").Value = -1097.3334
$ws.Range("N97").Value = -2995.3334
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1249.1428
$ws.Range("I7").Value = 1309.3334
$ws.Range("J7").Value = 888
$ws.Range("K7").Value = 1309.3334
$ws.Range("L7").Value = 888
$ws.Range("M7").Value = -1197.3334
$ws.Range("N7").Value = -1112
$ws.Range("H40").Value = 2269.2222
$ws.Range("I40").Value = 2346.1428
$ws.Range("K40").Value = 2346.1428
$ws.Range("M40").Value = -2210.1428
$ws.Range("H46").Value = 4718.6875
$ws.Range("I46").Value = 4000
$ws.Range("K46").Value = 4000
$ws.Range("M46").Value = -3812
$ws.Range("H61").Value = 849.75
$ws.Range("I61").Value = 833
$ws.Range("K61").Value = 833
$ws.Range("M61").Value = -631
$ws.Range("H100").Value = 12733
$ws.Range("J100").Value = 17850
$ws.Range("L100").Value = 17850
$ws.Range("N100").Value = -18932
$ws.Range("H113").Value = 849.75
$ws.Range("I113").Value = 833
$ws.Range("K113").Value = 833
$ws.Range("M113").Value = 1337
$ws.Range("H126").Value = 1249.1428
$ws.Range("I126").Value = 1309.3334
$ws.Range("J126").Value = 888
$ws.Range("K126").Value = 3928.0002
$ws.Range("L126").Value = 2664
$ws.Range("M126").Value = -1458.0002
$ws.Range("N126").Value = -7604
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H75").Value = 25000
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 25000
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H100").Value = 6338409
$ws.Range("I100").Value = 8713644
$ws.Range("J100").Value = 4450
$ws.Range("K100").Value = 17427288
$ws.Range("L100").Value = 8900
$ws.Range("M100").Value = -17426747
$ws.Range("N100").Value = -9982

Write-Output "applied all edits"